# Loop-style fill of names into the sheet, matching the evolving
# "teste1.xlsx" automation experiment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New client names appended below the existing list (A5:A7)
$ws.Range("A5").Value = "Amanda"
$ws.Range("A6").Value = "Thais"
$ws.Range("A7").Value = "Bernadete"

# One more write further down/over in column G, reusing "Thais"
$ws.Range("G19").Value = "Thais"

# Select A8:A11 (active cell ends on A11 once the loop finishes)
$ws.Range("A8:A11").Select()

# Switch the page to portrait orientation
$ws.PageSetup.Orientation = 1
